$wb = $excel.ActiveWorkbook

# ===== Sheet: Active Signals =====
$ws1 = $wb.Worksheets.Item("Active Signals")

# Create row 9 by copying formatting from row 2 (SELL style on column C matches needed style for new row)
$ws1.Range("A2:J2").Copy()
$ws1.Range("A9:J9").PasteSpecial(-4122)

# --- Row 2 ---
$ws1.Cells(2,1).Value = "2025-07-28 19:22"
$ws1.Cells(2,2).Value = "XAUCHF"
$ws1.Cells(2,3).Value = "SELL"
$ws1.Cells(2,4).Value = 2334.21327
$ws1.Cells(2,5).Value = 2334.21685
$ws1.Cells(2,6).Value = 2334.20676
$ws1.Cells(2,7).Value = 0.07
$ws1.Cells(2,8).Formula = '="93.0%"'
$ws1.Cells(2,9).Value = 1.82
$ws1.Cells(2,10).Value = "Active"

# --- Row 3 ---
$ws1.Cells(3,1).Value = "2025-07-28 19:57"
$ws1.Cells(3,2).Value = "USDJPY"
$ws1.Cells(3,3).Value = "BUY"
$ws1.Cells(3,4).Value = 149.85394
$ws1.Cells(3,5).Value = 149.37612
$ws1.Cells(3,6).Value = 150.65947
$ws1.Cells(3,7).Value = 0.04
$ws1.Cells(3,8).Formula = '="86.0%"'
$ws1.Cells(3,9).Value = 1.69
$ws1.Cells(3,10).Value = "Active"

# --- Row 4 ---
$ws1.Cells(4,1).Value = "2025-07-28 20:14"
$ws1.Cells(4,2).Value = "XAUCHF"
$ws1.Cells(4,3).Value = "BUY"
$ws1.Cells(4,4).Value = 2330.99258
$ws1.Cells(4,5).Value = 2330.98966
$ws1.Cells(4,6).Value = 2331.00069
$ws1.Cells(4,7).Value = 0.08
$ws1.Cells(4,8).Formula = '="87.0%"'
$ws1.Cells(4,9).Value = 2.78
$ws1.Cells(4,10).Value = "Active"

# --- Row 5 ---
$ws1.Range("C4").Copy()
$ws1.Range("C5").PasteSpecial(-4122)
$ws1.Cells(5,1).Value = "2025-07-28 20:08"
$ws1.Cells(5,2).Value = "XAUEUR"
$ws1.Cells(5,3).Value = "BUY"
$ws1.Cells(5,4).Value = 2425.36463
$ws1.Cells(5,5).Value = 2425.36158
$ws1.Cells(5,6).Value = 2425.37288
$ws1.Cells(5,7).Value = 0.01
$ws1.Cells(5,8).Formula = '="85.0%"'
$ws1.Cells(5,9).Value = 2.7
$ws1.Cells(5,10).Value = "Active"

# --- Row 6 ---
$ws1.Cells(6,1).Value = "2025-07-28 19:55"
$ws1.Cells(6,2).Value = "USDJPY"
$ws1.Cells(6,3).Value = "SELL"
$ws1.Cells(6,4).Value = 148.85479
$ws1.Cells(6,5).Value = 149.12514
$ws1.Cells(6,6).Value = 147.95812
$ws1.Cells(6,7).Value = 0.08
$ws1.Cells(6,8).Formula = '="88.0%"'
$ws1.Cells(6,9).Value = 3.32
$ws1.Cells(6,10).Value = "Active"

# --- Row 7 ---
$ws1.Range("C6").Copy()
$ws1.Range("C7").PasteSpecial(-4122)
$ws1.Cells(7,1).Value = "2025-07-28 20:19"
$ws1.Cells(7,2).Value = "USDCHF"
$ws1.Cells(7,3).Value = "SELL"
$ws1.Cells(7,4).Value = 0.88023
$ws1.Cells(7,5).Value = 0.88297
$ws1.Cells(7,6).Value = 0.87613
$ws1.Cells(7,7).Value = 0.03
$ws1.Cells(7,8).Formula = '="90.0%"'
$ws1.Cells(7,9).Value = 1.49
$ws1.Cells(7,10).Value = "Active"

# --- Row 8 ---
$ws1.Cells(8,1).Value = "2025-07-28 20:01"
$ws1.Cells(8,2).Value = "AUDUSD"
$ws1.Cells(8,3).Value = "BUY"
$ws1.Cells(8,4).Value = 0.66044
$ws1.Cells(8,5).Value = 0.65635
$ws1.Cells(8,6).Value = 0.66739
$ws1.Cells(8,7).Value = 0.03
$ws1.Cells(8,8).Formula = '="81.0%"'
$ws1.Cells(8,9).Value = 1.7
$ws1.Cells(8,10).Value = "Active"

# --- Row 9 ---
$ws1.Cells(9,1).Value = "2025-07-28 19:56"
$ws1.Cells(9,2).Value = "XAUUSD"
$ws1.Cells(9,3).Value = "SELL"
$ws1.Cells(9,4).Value = 2657.19974
$ws1.Cells(9,5).Value = 2657.20194
$ws1.Cells(9,6).Value = 2657.19485
$ws1.Cells(9,7).Value = 0.03
$ws1.Cells(9,8).Formula = '="78.0%"'
$ws1.Cells(9,9).Value = 2.23
$ws1.Cells(9,10).Value = "Active"

# Convert H-column formulas (percent-as-text trick) into plain text values, preserving style
$ws1.Range("H2:H9").Copy()
$ws1.Range("H2:H9").PasteSpecial(-4163)

# ===== Sheet: Summary Dashboard =====
$ws2 = $wb.Worksheets.Item("Summary Dashboard")
$ws2.Cells(4,2).Value = 8
$ws2.Cells(7,2).Formula = '="81.1%"'
$ws2.Cells(8,2).Formula = '="2.15"'
$ws2.Cells(9,2).Value = "2025-07-28 19:51:42"
$ws2.Range("B7:B8").Copy()
$ws2.Range("B7:B8").PasteSpecial(-4163)

# ===== Sheet: Signal History =====
$ws3 = $wb.Worksheets.Item("Signal History")

# --- Row 2 ---
$ws3.Cells(2,1).Value = "2025-07-28 20:07"
$ws3.Cells(2,2).Value = "XAUUSD"
$ws3.Cells(2,3).Value = "BUY"
$ws3.Cells(2,4).Value = 2639.10754
$ws3.Cells(2,5).Value = 2639.10473
$ws3.Cells(2,6).Value = 2639.11523
$ws3.Cells(2,7).Value = 0.03
$ws3.Cells(2,8).Value = 0.76
$ws3.Cells(2,9).Value = 2.73
$ws3.Cells(2,10).Value = "Filled"

# --- Row 3 ---
$ws3.Cells(3,1).Value = "2025-07-28 19:22"
$ws3.Cells(3,2).Value = "XAUCHF"
$ws3.Cells(3,3).Value = "SELL"
$ws3.Cells(3,4).Value = 2334.21327
$ws3.Cells(3,5).Value = 2334.21685
$ws3.Cells(3,6).Value = 2334.20676
$ws3.Cells(3,7).Value = 0.07
$ws3.Cells(3,8).Value = 0.93
$ws3.Cells(3,9).Value = 1.82
$ws3.Cells(3,10).Value = "Active"

# --- Row 4 ---
$ws3.Cells(4,1).Value = "2025-07-28 19:43"
$ws3.Cells(4,2).Value = "USDJPY"
$ws3.Cells(4,3).Value = "SELL"
$ws3.Cells(4,4).Value = 148.7955
$ws3.Cells(4,5).Value = 149.02496
$ws3.Cells(4,6).Value = 148.32316
$ws3.Cells(4,7).Value = 0.06
$ws3.Cells(4,8).Value = 0.88
$ws3.Cells(4,9).Value = 2.06
$ws3.Cells(4,10).Value = "Filled"

# --- Row 5 ---
$ws3.Cells(5,1).Value = "2025-07-28 19:57"
$ws3.Cells(5,2).Value = "USDJPY"
$ws3.Cells(5,3).Value = "BUY"
$ws3.Cells(5,4).Value = 149.85394
$ws3.Cells(5,5).Value = 149.37612
$ws3.Cells(5,6).Value = 150.65947
$ws3.Cells(5,7).Value = 0.04
$ws3.Cells(5,8).Value = 0.86
$ws3.Cells(5,9).Value = 1.69
$ws3.Cells(5,10).Value = "Active"

# --- Row 6 ---
$ws3.Cells(6,1).Value = "2025-07-28 20:14"
$ws3.Cells(6,2).Value = "XAUCHF"
$ws3.Cells(6,3).Value = "BUY"
$ws3.Cells(6,4).Value = 2330.99258
$ws3.Cells(6,5).Value = 2330.98966
$ws3.Cells(6,6).Value = 2331.00069
$ws3.Cells(6,7).Value = 0.08
$ws3.Cells(6,8).Value = 0.87
$ws3.Cells(6,9).Value = 2.78
$ws3.Cells(6,10).Value = "Active"

# --- Row 7 ---
$ws3.Cells(7,1).Value = "2025-07-28 20:04"
$ws3.Cells(7,2).Value = "NZDUSD"
$ws3.Cells(7,3).Value = "BUY"
$ws3.Cells(7,4).Value = 0.58889
$ws3.Cells(7,5).Value = 0.59119
$ws3.Cells(7,6).Value = 0.58034
$ws3.Cells(7,7).Value = 0.08
$ws3.Cells(7,8).Value = 0.73
$ws3.Cells(7,9).Value = 3.72
$ws3.Cells(7,10).Value = "Filled"

# --- Row 8 ---
$ws3.Cells(8,1).Value = "2025-07-28 20:08"
$ws3.Cells(8,2).Value = "XAUEUR"
$ws3.Cells(8,3).Value = "BUY"
$ws3.Cells(8,4).Value = 2425.36463
$ws3.Cells(8,5).Value = 2425.36158
$ws3.Cells(8,6).Value = 2425.37288
$ws3.Cells(8,7).Value = 0.01
$ws3.Cells(8,8).Value = 0.85
$ws3.Cells(8,9).Value = 2.7
$ws3.Cells(8,10).Value = "Active"

# --- Row 9 ---
$ws3.Cells(9,1).Value = "2025-07-28 19:39"
$ws3.Cells(9,2).Value = "XAUAUD"
$ws3.Cells(9,3).Value = "BUY"
$ws3.Cells(9,4).Value = 4032.93353
$ws3.Cells(9,5).Value = 4032.9298
$ws3.Cells(9,6).Value = 4032.93825
$ws3.Cells(9,7).Value = 0.03
$ws3.Cells(9,8).Value = 0.65
$ws3.Cells(9,9).Value = 1.26
$ws3.Cells(9,10).Value = "Pending"

# --- Row 10 ---
$ws3.Cells(10,1).Value = "2025-07-28 19:47"
$ws3.Cells(10,2).Value = "USDJPY"
$ws3.Cells(10,3).Value = "SELL"
$ws3.Cells(10,4).Value = 149.13696
$ws3.Cells(10,5).Value = 149.62288
$ws3.Cells(10,6).Value = 148.65307
$ws3.Cells(10,7).Value = 0.07
$ws3.Cells(10,8).Value = 0.68
$ws3.Cells(10,9).Value = 1
$ws3.Cells(10,10).Value = "Pending"

# --- Row 11 ---
$ws3.Cells(11,1).Value = "2025-07-28 19:55"
$ws3.Cells(11,2).Value = "USDJPY"
$ws3.Cells(11,3).Value = "SELL"
$ws3.Cells(11,4).Value = 148.85479
$ws3.Cells(11,5).Value = 149.12514
$ws3.Cells(11,6).Value = 147.95812
$ws3.Cells(11,7).Value = 0.08
$ws3.Cells(11,8).Value = 0.88
$ws3.Cells(11,9).Value = 3.32
$ws3.Cells(11,10).Value = "Active"

# --- Row 12 ---
$ws3.Cells(12,1).Value = "2025-07-28 20:19"
$ws3.Cells(12,2).Value = "USDCHF"
$ws3.Cells(12,3).Value = "SELL"
$ws3.Cells(12,4).Value = 0.88023
$ws3.Cells(12,5).Value = 0.88297
$ws3.Cells(12,6).Value = 0.87613
$ws3.Cells(12,7).Value = 0.03
$ws3.Cells(12,8).Value = 0.9
$ws3.Cells(12,9).Value = 1.49
$ws3.Cells(12,10).Value = "Active"

# --- Row 13 ---
$ws3.Cells(13,1).Value = "2025-07-28 20:01"
$ws3.Cells(13,2).Value = "AUDUSD"
$ws3.Cells(13,3).Value = "BUY"
$ws3.Cells(13,4).Value = 0.66044
$ws3.Cells(13,5).Value = 0.65635
$ws3.Cells(13,6).Value = 0.66739
$ws3.Cells(13,7).Value = 0.03
$ws3.Cells(13,8).Value = 0.81
$ws3.Cells(13,9).Value = 1.7
$ws3.Cells(13,10).Value = "Active"

# --- Row 14 ---
$ws3.Cells(14,1).Value = "2025-07-28 20:01"
$ws3.Cells(14,2).Value = "XAUEUR"
$ws3.Cells(14,3).Value = "SELL"
$ws3.Cells(14,4).Value = 2412.71651
$ws3.Cells(14,5).Value = 2412.71944
$ws3.Cells(14,6).Value = 2412.71005
$ws3.Cells(14,7).Value = 0.06
$ws3.Cells(14,8).Value = 0.92
$ws3.Cells(14,9).Value = 2.21
$ws3.Cells(14,10).Value = "Filled"

# --- Row 15 ---
$ws3.Cells(15,1).Value = "2025-07-28 19:25"
$ws3.Cells(15,2).Value = "GBPUSD"
$ws3.Cells(15,3).Value = "SELL"
$ws3.Cells(15,4).Value = 1.26529
$ws3.Cells(15,5).Value = 1.26994
$ws3.Cells(15,6).Value = 1.25782
$ws3.Cells(15,7).Value = 0.04
$ws3.Cells(15,8).Value = 0.67
$ws3.Cells(15,9).Value = 1.6
$ws3.Cells(15,10).Value = "Filled"

# --- Row 16 ---
$ws3.Cells(16,1).Value = "2025-07-28 19:56"
$ws3.Cells(16,2).Value = "XAUUSD"
$ws3.Cells(16,3).Value = "SELL"
$ws3.Cells(16,4).Value = 2657.19974
$ws3.Cells(16,5).Value = 2657.20194
$ws3.Cells(16,6).Value = 2657.19485
$ws3.Cells(16,7).Value = 0.03
$ws3.Cells(16,8).Value = 0.78
$ws3.Cells(16,9).Value = 2.23
$ws3.Cells(16,10).Value = "Active"
